# TweetDeckSPFxFilter.xlsx edit script
# - Renames the sheet Tabelle1 -> SPFx
# - Adds a label/formula/help block at the top of the sheet (rows 1-4)
# - Adds a LEN() character-count formula with conditional formatting (500 char Tweetdeck limit)
# - Updates the Hashtags list (adds/removes several hashtags, keeps it alphabetically sorted)
# - Resizes the Hashtags/From tables to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: push the existing table down by 3 rows (for the new header
#    block) and remember the style of the two cells that used to live at the
#    bottom of the sheet (the "TweetDeckSPFxFilter" label and the big
#    CONCATENATE formula) so the new header cells can reuse them exactly.
# ---------------------------------------------------------------------------
$ws.Range("A1:E3").EntireRow.Insert()

# After the insert above, the old E30/E31 (label + formula) now live at
# E33/E34. Copy their formatting onto the new E1/E2 before removing them.
$ws.Range("E33").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E34").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("A33:E34").EntireRow.Delete()

# Make room for the extra hashtags (28 -> 37 entries) below the current list.
$ws.Range("A33:A41").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Resize the two structured tables to their final ranges before writing
#    any formulas that use structured references, so the references resolve
#    against the correct ranges.
# ---------------------------------------------------------------------------
$loHashtags = $ws.ListObjects.Item("Hashtags")
$loFrom = $ws.ListObjects.Item("From")
$loHashtags.Resize($ws.Range("A4:A41"))
$loFrom.Resize($ws.Range("C4:C6"))

# ---------------------------------------------------------------------------
# 3. Rewrite the Hashtags column with the updated, alphabetically sorted list.
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value = "Hashtags"

$hashtags = @(
  "#art",
  "#artist",
  "#artistreborn",
  "#behindthescenes",
  "#blood",
  "#Elementary",
  "#fxmakeup",
  "#gore",
  "#gross",
  "#hollywood",
  "#horror",
  "#LED",
  "#makeup",
  "#makeupartist",
  "#makeupeffects",
  "#makeupfx",
  "#mua",
  "#mufx",
  "#puppets",
  "#pyro",
  "#pyrotechnics",
  "#sculpture",
  "#sfx",
  "#sfxmakeup",
  "#SharpFX",
  "#specialeffects",
  "#specialfx",
  "#specialfxmakeup",
  "#specialmakeupfx",
  "#spfmakeup",
  "#spfxmakeup",
  "#spfxmakeupartist",
  "#spfxmua",
  "#thrillerthursday",
  "#vfx",
  "#werewolf",
  "#wounds"
)

for ($i = 0; $i -lt $hashtags.Length; $i++) {
  $ws.Cells.Item(5 + $i, 1).Value = $hashtags[$i]
}

# ---------------------------------------------------------------------------
# 4. "From" table header/data (unchanged content, just confirm placement).
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 3).Value = "From"
$ws.Cells.Item(5, 3).Value = "webpart_o_matic"

# ---------------------------------------------------------------------------
# 5. Header block in column E: label, filter formula, help text, length check.
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "TweetDeckSPFxFilter"

$ws.Range("E2").Formula = '=CONCATENATE("#Office365Dev OR @OfficeDev OR from:OfficeDev OR #SPFx", IF(COUNTIF(Hashtags[],"*")>0," -",""),_xlfn.TEXTJOIN(" -", TRUE,Hashtags[]), IF(COUNTIF(From[],"*")>0," -from:",""),_xlfn.TEXTJOIN(" -from:", TRUE, From[]))'

$ws.Range("E3").Value = "Max lenght of Tweetdeck filter is 500 characters. Current filter lenght:"

$ws.Range("E4").Formula = "=LEN(E2)"

# ---------------------------------------------------------------------------
# 6. Conditional formatting on E4: highlight red if over the 500 char limit,
#    green if comfortably under it (classic Highlight-Cell-Rules colours).
# ---------------------------------------------------------------------------
$cfRange = $ws.Range("E4")
$cfGreater = $cfRange.FormatConditions.Add(1, 5, "500")
$cfGreater.Interior.Color = 13551615
$cfGreater.Font.Color = 393372
$cfLess = $cfRange.FormatConditions.Add(1, 6, "500")
$cfLess.Interior.Color = 13561798
$cfLess.Font.Color = 24832
$cfLess.SetFirstPriority()

# ---------------------------------------------------------------------------
# 7. Cosmetics: column widths and sheet/tab name.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.833333
$ws.Columns.Item(2).ColumnWidth = 0.5
$ws.Columns.Item(3).ColumnWidth = 14.666667
$ws.Columns.Item(4).ColumnWidth = 0.5
$ws.Columns.Item(5).ColumnWidth = 101.833333
$ws.Columns.Item(8).ColumnWidth = 17.333333

[void]$ws.Range("C6").Select()

$ws.Name = "SPFx"

Write-Host "Edit complete"
